$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.099.26'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '2.472.74'
$ws.Range("E3").Value = '  +2.35%  '
$ws.Range("E4").Value = '  -0.47%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.79'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.06%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -0.40%  '
$ws.Range("D9").Value = '2.471.79'
$ws.Range("E9").Value = '  +1.03%  '
$ws.Range("E10").Value = '  +0.54%  '
$ws.Range("E11").Value = '  +1.71%  '
$ws.Range("E12").Value = '  +1.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.352'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000179'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = '2.921.63'
$ws.Range("E16").Value = '  +1.89%  '
$ws.Range("D17").Value = '62.985.95'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '2.474.19'
$ws.Range("E18").Value = '  +2.15%  '
$ws.Range("E19").Value = '  +3.81%  '
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '329.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("E22").Value = '  +10.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '663.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.61'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.27%  '
$ws.Range("D28").Value = '0.0₃0989'
$ws.Range("E28").Value = '  +0.90%  '
$ws.Range("D29").Value = '2.593.49'
$ws.Range("E29").Value = '  +2.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.65%  '
$ws.Range("E31").Value = '  +2.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.06'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("E33").Value = '  +1.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.133'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.54'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.81%  '
$ws.Range("E36").Value = '  +0.34%  '
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("E38").Value = '  +0.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '152.48'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.25%  '
$ws.Range("E40").Value = '  -0.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.78'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.75'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").Value = '0.0₆0306'
$ws.Range("E45").Value = '  +2.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '151.88'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.81%  '
$ws.Range("E47").Value = '  +24.97%  '
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.606'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0512'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.76%  '
